$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 5 (including all of its cell formatting/styles) down into a new
# row 6 so the new BOM line picks up the exact same styles (currency format,
# vertically-centered text format, etc.) already used by the other data rows.
$ws.Rows(5).Copy()
$ws.Rows(6).Insert(-4121, -4163)
$excel.CutCopyMode = 0

# Populate the new row for the heatshrink part. H6 (the link) is written first
# so the shared-string table ends up with the same ordering as the source
# workbook (link, then description, manufacturer, part number).
$ws.Cells.Item(6, 8).Value = "https://www.mouser.com/ProductDetail/CUI/HSE-B20250-040H?qs=sGAEpiMZZMttgyDkZ5Wiut%252B4GcHIZ2pKOgousR6bMSo%3D"
$ws.Cells.Item(6, 2).Value = "Heatshrink"
$ws.Cells.Item(6, 3).Value = "CUI"
$ws.Cells.Item(6, 4).Value = "HSE-B20250-040H "
$ws.Cells.Item(6, 5).Value = 0.58
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(6, 7).Formula = "=F6*E6"

# Unlike the other "link" cells in column H, this one is plain text (no
# hyperlink relationship), so strip the hyperlink formatting that got copied
# down from row 5.
$ws.Cells.Item(6, 8).Style = "Normal"

# Match the saved selection/active cell of the edited workbook.
$ws.Range("B2:H2").Select()
